$d = $word.ActiveDocument

# --- Color updates (RGB 0070C0 -> OLE BGR decimal 12611584) ---
$blueColor = 12611584

# Paragraph: "สามารถดูคนที่มา Follow เราได้ และ Tag อะไรบ้างที่เราติดตาม (หน้า Profile)"
$d.Paragraphs.Item(11).Range.Font.Color = $blueColor

# Paragraph: "หน้า Profile แสดงยอด Like ของบทความทั้งหมด ตำแหน่งถัดจาก Follows"
$d.Paragraphs.Item(12).Range.Font.Color = $blueColor

# Paragraph: "มีระบบ Suggested (3 บทความ มาจาก Tag ที่เหมือนกัน)"
$d.Paragraphs.Item(17).Range.Font.Color = $blueColor

# Paragraph: "สามารถดูรายชื่อคนที่มา Follow เราได้"
$d.Paragraphs.Item(21).Range.Font.Color = $blueColor

# --- Move the "_GoBack" bookmark from the second "Article" paragraph ---
# to just after the text of the first "Article" paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$articlePara = $d.Paragraphs.Item(13)
$insertPos = $articlePara.Range.End - 1
$tmp = $d.Range($insertPos, $insertPos)
$tmp.InsertAfter("X")
$markRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$bmRange = $d.Bookmarks.Item("_GoBack").Range
$bmRange.Text = ""

Write-Host "All edits applied"
